$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cells in row 2 (D2, E2)
$ws.Range("D2").Value = "Down 6"
$ws.Range("E2").Value = "Come back to 3.5 Automated testing with Mocha"

# Widen column E like the other data columns
$ws.Columns.Item(5).ColumnWidth = 22.85

# New rows appended after the existing data (row 25 was the last one before)
$ws.Range("A26").Value = "3.5 Automated testing with Mocha skipped"
$ws.Range("A27").Value = "3.6 Polyfills"
# Row 28 intentionally left blank
$ws.Range("A29").Value = "4.2 Garbage collection"

# Distinguish A26 with its own style (new font/style entry observed in the diff)
$ws.Range("A26").Font.Bold = $true

# Update selection to match diff (A29 active cell)
$ws.Range("A29").Select()
